$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values on both the "展览" and "全部类型"
# sheets, which contain duplicate data.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 134
    $ws.Range("F3").Value = 22
    $ws.Range("F5").Value = 26
}
